# BRE 2.0.1 update — re-parametrize "PerfilClienteExterno" score bands
# and refresh the sheet view state left behind by the editing session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ---------------------------------------------------------------------
# 1) Re-parametrize the Equifax score bands (columns C/D) and the
#    resulting GRUPO values (column F) for the external-client profile
#    rule table.
# ---------------------------------------------------------------------
$ws.Range("C9").Value = 818

$ws.Range("C10").Value = 656
$ws.Range("D10").Value = 817

$ws.Range("C11").Value = 468
$ws.Range("D11").Value = 655
$ws.Range("F11").Value = 2

$ws.Range("C12").Value = 231
$ws.Range("D12").Value = 467
$ws.Range("F12").Value = 3

$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 230
$ws.Range("F13").Value = 3

$ws.Range("C14").Value = -1
$ws.Range("D14").Value = -1
$ws.Range("F14").Value = 3

$ws.Range("C15").Value = -1
$ws.Range("D15").Value = -1
$ws.Range("F15").Value = 3

$ws.Range("C16").Value = -1
$ws.Range("D16").Value = -1
$ws.Range("F16").Value = 3

$ws.Range("C17").Value = -1
$ws.Range("D17").Value = -1

$ws.Range("C18").Value = -1
$ws.Range("D18").Value = -1

# ---------------------------------------------------------------------
# 2) Fix stray header-cell formatting: D5 and E7 had picked up a
#    duplicated/odd style at some point; re-apply the same left-aligned
#    "CONDITION/ACTION" header look used by their row neighbours.
# ---------------------------------------------------------------------
$ws.Range("C5").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("C7").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Restore the working view: zoomed out to 70% and the cursor left on
#    the last edited parametrization row.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 70
$ws.Range("D17").Select() | Out-Null
